$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update task start dates (column C) ---
# Row 16 ("Stworzenie GUI") slipped by one day
$ws.Range("C16").Value = 42834
# Row 18 ("Testy całości") slipped by one day
$ws.Range("C18").Value = 42835
# Row 19 ("Napisanie dokumentacji") slipped by one day
$ws.Range("C19").Value = 42836
# Row 20 ("Napisanie Post Mortem") slipped by one day
$ws.Range("C20").Value = 42836

# --- Mark rows 17 and 18 as done (column E), matching the style/marker
# already used for the other completed tasks (e.g. row 16) ---
$ws.Range("E16").Copy()
$ws.Range("E17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("E17").Value = "þ"
$ws.Range("E18").Value = "þ"

# --- Update the active selection to reflect where the user left off ---
$ws.Range("F17").Select()
